# Auto-generated edit script applying the cryptos.xlsx cell-value update
# (price/volume refresh + a few re-ranked rows) described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.380.95'
$ws.Range('E2').Value = '  +2.49%  '
$ws.Range('D3').Value = '3.784.58'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'624.75"
$ws.Range('E5').Value = '  +4.38%  '
$ws.Range('D6').Value = "'165.77"
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('D7').Value = '3.782.40'
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.519"
$ws.Range('E9').Value = '  +1.32%  '
$ws.Range('E10').Value = '  +2.80%  '
$ws.Range('D11').Value = "'0.453"
$ws.Range('E11').Value = '  +1.98%  '
$ws.Range('D12').Value = "'6.70"
$ws.Range('E12').Value = '  +2.38%  '
$ws.Range('E13').Value = '  +1.43%  '
$ws.Range('D14').Value = "'35.68"
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').Value = '4.420.19'
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').Value = '3.795.62'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').Value = '69.362.40'
$ws.Range('E17').Value = '  +2.44%  '
$ws.Range('D18').Value = "'17.70"
$ws.Range('E18').Value = '  -2.92%  '
$ws.Range('E19').Value = '  +2.00%  '
$ws.Range('D20').Value = "'0.114"
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('D21').Value = "'468.29"
$ws.Range('E21').Value = '  +2.48%  '
$ws.Range('D22').Value = "'9.63"
$ws.Range('E22').Value = '  +1.70%  '
$ws.Range('D23').Value = "'0.705"
$ws.Range('E23').Value = '  +1.82%  '
$ws.Range('E24').Value = '  +4.71%  '
$ws.Range('D25').Value = "'83.30"
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('D26').Value = "'12.03"
$ws.Range('E26').Value = '  +1.50%  '
$ws.Range('E27').Value = '  +3.89%  '
$ws.Range('E28').Value = '  +1.93%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').Value = '3.933.77'
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('E31').Value = '  +3.21%  '
$ws.Range('E33').Value = '  +0.72%  '
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = "'0.167"
$ws.Range('E35').Value = '  +18.42%  '
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').Value = "'1.00"
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').Value = '3.735.35'
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = "'8.98"
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('E39').Value = '  +3.12%  '
$ws.Range('D40').Value = "'3.41"
$ws.Range('E40').Value = '  +8.45%  '
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('D42').Value = "'0.969"
$ws.Range('E42').Value = '  -0.85%  '
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D45').Value = "'43.45"
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = "'1.92"
$ws.Range('E47').Value = '  +4.50%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = "'46.78"
$ws.Range('E48').Value = '  -0.84%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = "'151.89"
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('E50').Value = '  +2.16%  '
$ws.Range('E51').Value = '  +0.68%  '

Write-Host "Updated 93 cells"
